# Add Inflation_contributions_graph_data.XLSX, Waterfall_graph_data.XLSX,
# trimmed_graph_data.XLSX, MEGA_DATA_DOWNLOAD.xlsx, MONTHLY_DATA_DOWNLOAD_ALL.xlsx,
# QUARTERLY_DATA_DOWNLOAD_ALL.xlsx, SUMMARY_EXPORT_DATA_DOWNLOAD.xlsx and
# EXPORT_DATA_DOWNLOAD_ALL.xlsx -- update the index-values worksheet with three
# additional monthly observations (and correct the most recent value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch cell, far away from the used range, that we use as a bounce pad so a
# numeric-looking piece of text (e.g. "110.8") is forced to be written as a real
# shared string rather than silently re-interpreted as a number.
$scratchRow = 5000
$scratchCol = 50

function Set-TextValue {
    param($targetRow, $targetCol, [string]$text)

    $scratch = $ws.Cells.Item($scratchRow, $scratchCol)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Cells.Item($targetRow, $targetCol).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

function Set-DateValue {
    param($targetRow, $targetCol, $serial)

    # Reuse the existing date formatting (style index used by column A) by
    # copying formats from the cell directly above, then set the numeric value.
    $src = $ws.Cells.Item($targetRow - 1, $targetCol)
    $src.Copy()
    $dst = $ws.Cells.Item($targetRow, $targetCol)
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $dst.Value = $serial
}

# --- Row 97: correct the most recent Services value -------------------------
$ws.Cells.Item(97, 2).Value = 130.912633615201

# --- Row 98: new monthly observation (2025-01-01) ----------------------------
Set-DateValue 98 1 45658
$ws.Cells.Item(98, 2).Value = 129.599960964251
$ws.Cells.Item(98, 3).Value = 122.427773483352
Set-TextValue 98 4 "110.8"
Set-TextValue 98 5 "112.9"
Set-TextValue 98 6 " 86.4"
Set-TextValue 98 7 "171.1"

# --- Row 99: new monthly observation (2025-02-01) ----------------------------
Set-DateValue 99 1 45689
$ws.Cells.Item(99, 2).Value = 130.763998761473
$ws.Cells.Item(99, 3).Value = 122.637687759991
Set-TextValue 99 4 "111.1"
Set-TextValue 99 5 "113.5"
Set-TextValue 99 6 " 87.5"
Set-TextValue 99 7 "171.3"

# --- Row 100: new monthly observation (2025-03-01) ---------------------------
Set-DateValue 100 1 45717
$ws.Cells.Item(100, 2).Value = 132.015597494155
$ws.Cells.Item(100, 3).Value = 123.167012106947
Set-TextValue 100 4 "112.1"
Set-TextValue 100 5 "114.3"
Set-TextValue 100 6 " 88.3"
Set-TextValue 100 7 "172.0"
